$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Table2")

$ws.Range("D2").Value = 3.7000000476837158
$ws.Range("D3").Value = 4.8000001907348633
$ws.Range("D4").Value = 5.5
$ws.Range("D5").Value = 5.8000001907348633
